$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 501, pushing existing rows 501-563 down to 502-564.
$ws.Rows.Item(501).Insert()

# Populate the newly inserted row 501 with the new weekly data point.
$ws.Range("A501").Value = 9
$ws.Range("B501").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C501").Value = "Metropolitana"
$ws.Range("D501").Value = 45077
$ws.Range("E501").Value = 13
$ws.Range("F501").Value = 100112044
$ws.Range("G501").Value = "Perejil"
$ws.Range("H501").Value = "Sin especificar"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 70
$ws.Range("K501").Value = 11000
$ws.Range("L501").Value = 13000
$ws.Range("M501").Value = 12000
$ws.Range("N501").Value = "$/docena de atados"
$ws.Range("O501").Value = "Región Metropolitana"
$ws.Range("P501").Value = 4000
$ws.Range("Q501").Value = 3
$ws.Range("R501").Value = "Hortaliza"
